$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": PORCELANATO column (M) for two clients.
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M19").Value = 411.24
$wsVentasGrupo.Range("M33").Value = 5238.25

# Sheet "VENTA MENSUAL": septiembre column (F) for two clients + total row.
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F19").Value = 740.24
$wsVentaMensual.Range("F33").Value = 5238.25
$wsVentaMensual.Range("F34").Value = 8565.690000000001

# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row (12) + TOTAL row (15).
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 6845.65
$wsCumplimiento.Range("E12").Value = 15588.1053751766
$wsCumplimiento.Range("F12").Value = 0.3051495340621771
$wsCumplimiento.Range("D15").Value = 8565.689999999999
$wsCumplimiento.Range("E15").Value = 30177.32881339592
$wsCumplimiento.Range("F15").Value = 0.2210898959953604
$wsCumplimiento.Columns.Item(6).ColumnWidth = 23.17
